# Apply the Oct 20 2023 cryptos-list refresh (GitHub Actions bot update).
# Updates Price (D) / Volume(1h) (E) figures, and for rows 48-51 also the
# Coin name (B) and Link (C) since four coins got reordered in the feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.406.95'
$ws.Range("E2").Value = '  +2.76%  '

$ws.Range("D3").Value = '1.604.98'
$ws.Range("E3").Value = '  +2.53%  '

$ws.Range("D4").Value = "'" + '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = "'" + '212.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.76%  '

$ws.Range("E6").Value = '  +6.66%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").Value = '  +6.68%  '

$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("E10").Value = '  +2.54%  '

$ws.Range("E11").Value = '  +2.44%  '

$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").Value = '1.833.20'
$ws.Range("E13").Value = '  +2.43%  '

$ws.Range("D14").Value = '1.605.52'
$ws.Range("E14").Value = '  +2.51%  '

$ws.Range("D15").Value = '29.433.49'
$ws.Range("E15").Value = '  +2.76%  '

$ws.Range("E17").Value = '  +1.76%  '

$ws.Range("D18").Value = "'" + '63.15'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.96%  '

$ws.Range("D19").Value = "'" + '240.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.29%  '

$ws.Range("D20").Value = "'" + '7.64'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.61%  '

$ws.Range("D22").Value = "'" + '0.998'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("E23").Value = '  +1.91%  '

$ws.Range("E24").Value = '  +1.90%  '

$ws.Range("D25").Value = "'" + '2.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").Value = "'" + '154.65'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.71%  '

$ws.Range("E27").Value = '  +4.83%  '

$ws.Range("D28").Value = "'" + '15.26'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.25%  '

$ws.Range("D29").Value = "'" + '6.37'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.39%  '

$ws.Range("D30").Value = "'" + '0.998'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("E31").Value = '  +2.37%  '

$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("E33").Value = '  +1.56%  '

$ws.Range("E34").Value = '  +3.85%  '

$ws.Range("D35").Value = '1.411.86'
$ws.Range("E35").Value = '  +1.57%  '

$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("D37").Value = "'" + '1.53'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.87%  '

$ws.Range("D38").Value = "'" + '2.81'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.80%  '

$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("E41").Value = '  +3.26%  '

$ws.Range("D42").Value = "'" + '1.96'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.92%  '

$ws.Range("D43").Value = "'" + '0.0489'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.88%  '

$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("D45").Value = "'" + '0.795'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.89%  '

$ws.Range("D46").Value = "'" + '52.71'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +22.20%  '

$ws.Range("D47").Value = "'" + '65.55'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.89%  '

$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.745.06'
$ws.Range("E48").Value = '  +2.59%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = "'" + '5.25'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.52%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'" + '86.63'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.67%  '

$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = "'" + '0.847'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.71%  '
